$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.500139951705933
$ws.Range("B1").Value = 1.708860039710999
$ws.Range("C1").Value = 1.718041300773621
$ws.Range("D1").Value = 2.148695945739746
$ws.Range("E1").Value = 3.241360187530518
